$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 19): Day 10, Date text "19/6/2024", 1.25 hours,
# Description "Finalizat tot, mai ramane doar o ultima refactorizare"
# (set D19 before B19 so shared strings are appended in the same order as
# the target workbook: description string first, then the date string)
$ws.Range("A19").Value = 10
$ws.Range("D19").Value = "Finalizat tot, mai ramane doar o ultima refactorizare"
$ws.Range("C19").Value = 1.25
$ws.Range("B19").Value = "19/6/2024"

# Update the selected cell to C21 (matches author's final selection state)
$ws.Range("C21").Select()
